$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.575.96"
$ws.Range("E2").Value = "  +3.86%  "

$ws.Range("D3").Value = "2.075.22"
$ws.Range("E3").Value = "  +3.49%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.19"
$ws.Range("E5").Value = "  +2.70%  "

$ws.Range("E6").Value = "  +1.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.42"
$ws.Range("E7").Value = "  +6.46%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +7.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.63"
$ws.Range("E10").Value = "  +0.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0828"
$ws.Range("E11").Value = "  +10.74%  "

$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.928"
$ws.Range("E13").Value = "  -2.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.78"
$ws.Range("E14").Value = "  +20.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.88"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").Value = "2.373.86"
$ws.Range("E16").Value = "  +3.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.68"
$ws.Range("E17").Value = "  +4.36%  "

$ws.Range("D18").Value = "2.069.67"
$ws.Range("E18").Value = "  +3.89%  "

$ws.Range("D19").Value = "37.459.28"
$ws.Range("E19").Value = "  +3.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.70"
$ws.Range("E20").Value = "  +2.06%  "

$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("E21").Value = "  +6.11%  "

$ws.Range("E22").Value = "  +4.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.21"
$ws.Range("E23").Value = "  +2.61%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("E25").Value = "  -1.85%  "

$ws.Range("E26").Value = "  +3.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("E27").Value = "  +3.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.73"
$ws.Range("E28").Value = "  +5.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.35"
$ws.Range("E29").Value = "  -2.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.130"
$ws.Range("E30").Value = "  +30.15%  "

$ws.Range("E31").Value = "  +2.39%  "

$ws.Range("E32").Value = "  +1.83%  "

$ws.Range("E33").Value = "  +4.58%  "

$ws.Range("E34").Value = "  +3.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.71"
$ws.Range("E35").Value = "  +4.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.54"
$ws.Range("E36").Value = "  +2.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.45"
$ws.Range("E37").Value = "  +11.44%  "

$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("E39").Value = "  +3.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.08"
$ws.Range("E40").Value = "  +31.68%  "

$ws.Range("E41").Value = "  +7.35%  "

$ws.Range("E42").Value = "  +3.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.02"
$ws.Range("E43").Value = "  +4.33%  "

$ws.Range("E44").Value = "  +4.95%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0220"
$ws.Range("E45").Value = "  +1.88%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.31"
$ws.Range("E46").Value = "  +1.84%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.01"
$ws.Range("E47").Value = "  +2.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.80"
$ws.Range("E48").Value = "  +1.47%  "

$ws.Range("D49").Value = "1.401.86"
$ws.Range("E49").Value = "  +2.36%  "

$ws.Range("E50").Value = "  +1.41%  "

$ws.Range("E51").Value = "  -0.83%  "
